$d = $word.ActiveDocument

function ReplaceFirst($oldText, $newText) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Text = $oldText
    $found = $find.Execute()
    if (-not $found) {
        Write-Host "NOT FOUND: $oldText"
        return $null
    }
    $rng = $find.Parent
    $rng.Text = $newText
    return $rng
}

# ---- Title ----
ReplaceFirst "Unveiling the Enigmatic Realm of Dark Matter" `
    "Reading Through the Mosaic of History: A Comprehensive Expedition of Humanity's Collective Past" | Out-Null

# ---- Author (merge "Dr" + "." + " Sophia Barnett" -> "Emma Carter") ----
ReplaceFirst "Dr. Sophia Barnett" "Emma Carter" | Out-Null

# ---- Email user part ----
ReplaceFirst "sophiab@darkmatterhub" "emmacarter@carters" | Out-Null

# ---- Body paragraph 1 ----
ReplaceFirst "In the cosmic tapestry, there exists a mysterious substance known as dark matter, an enigmatic entity shrouding the universe in questions" `
    "The annals of history, a mosaic spanning centuries, invite us on an enthralling expedition" | Out-Null

ReplaceFirst " Its presence, though unseen, is felt in its gravitational influence on visible matter, shaping the dynamics and structure of galaxies and clusters" `
    " We delve into realms where civilizations rise and fall, leaders leave indelible marks, and cultures intertwine in a symphony of human experiences" | Out-Null

$rng = ReplaceFirst " It is believed to constitute approximately 27% of the universe's mass-energy, dwarfing the contribution of ordinary matter" `
    " In this odyssey of historical exploration, we unlock the secrets of ancient empires, decipher the echoes of forgotten eras, and ponder the relentless tides of change"
# The run immediately following is a lone "." run; extend it with the extra sentence
# that was inserted in the edit (new runs: "." then " Through the tapestry...future")
$endPos = $rng.End
$periodRng = $d.Range($endPos, $endPos + 1)
$periodRng.Text = ". Through the tapestry of history, we discover ourselves and forge connections with our collective past, charting the course for our journey into the future."

ReplaceFirst "Dark matter remains an enigma, a testament to human limitations in understanding the universe's fundamental nature" `
    "We scrutinize the intricacies of historical contexts, examining the forces that mold nations and the individuals who shape destinies" | Out-Null

ReplaceFirst " While its existence has been extensively inferred through gravitational lensing, rotation curves of galaxies, and cosmic microwave background radiation, its composition and properties remain elusive, fueling fervent scientific inquiry" `
    " From indomitable warriors to pioneering innovators, from enlightened rulers to ardent revolutionaries, the actors on history's stage teach us about courage, resilience, and the power of human agency" | Out-Null

$rng2 = ReplaceFirst " The search for a comprehensive theory that unifies the enigmatic world of dark matter with the Standard Model of Physics is a major quest of modern cosmology" `
    " We learn from their triumphs and misfortunes, their brilliance and follies, and grasp the intricate web of cause and effect that weaves the narrative of the past"

ReplaceFirst "The pursuit of dark matter unravels like a thrilling detective story" `
    "The study of history isn't just an academic pursuit; it's a transformative experience that shapes our understanding of the world" | Out-Null

ReplaceFirst " Evidence, like clues scattered across the vast expanse of the universe, guides scientists in their relentless quest" `
    " By embarking on this historical voyage, we develop critical thinking skills, cultural sensitivity, and a profound appreciation for the richness of human existence" | Out-Null

$rng3 = ReplaceFirst " As our telescopes peer deeper into the cosmos, and particle accelerators probe the subatomic world with ever-increasing precision, the mystery of dark matter inches closer to resolution, promising a fundamental understanding of the universe's fabric" `
    " We become global citizens, capable of navigating the complexities of the present with a deeper awareness of our shared history and an informed vision for the future"
$endPos3 = $rng3.End
$periodRng3 = $d.Range($endPos3, $endPos3 + 1)
$periodRng3.Text = ". We develop empathy, understanding the struggles and aspirations of people across time and place."

# ---- Summary paragraph ----
ReplaceFirst "Dark matter, a cosmic puzzle cloaked in mystery, exerts its gravitational influence throughout the universe, influencing the behavior of visible matter" `
    "In this essay, we embarked on an exhilarating journey through the mosaic of history, exploring the rise and fall of civilizations, the impact of influential leaders, and the intricate connections between cultures" | Out-Null

ReplaceFirst " Its elusive nature, comprising approximately 27% of the universe's mass-energy, challenges our understanding of fundamental physics" `
    " We delved into the forces that shape societies, the narratives of triumph and tragedy, and the lessons we can glean from the past" | Out-Null

ReplaceFirst " Through observations, simulations, and experiments, scientists relentlessly pursue the truth behind dark matter, aiming to unravel its composition and properties" `
    " By understanding history, we deepen our appreciation for humanity's collective experience, cultivate critical thinking skills, and forge a connection with our shared heritage" | Out-Null

ReplaceFirst " Unveiling the secrets of dark matter promises a deeper comprehension of the universe's structure and evolution, opening new vistas in our perception of reality" `
    " History becomes a mirror through which we learn about ourselves, our world, and the boundless possibilities that lie ahead" | Out-Null

# ---- Append a new empty paragraph at the end of the document ----
$end = $d.Content.End
$tail = $d.Range($end - 1, $end - 1)
$tail.InsertParagraphAfter()

Write-Host "done"
